$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content/formatting so the sheet can be rebuilt to match the
# "Completed Tier 1 Backend" layout from a clean slate.
$ws.Range("A1:J50").Clear()

# --- Cell values ---
$ws.Range("A2").Value = "Tier 1 Goal - Dashboard"
$ws.Range("B2").Value = "Database"
$ws.Range("C2").Value = "Express"
$ws.Range("D2").Value = "React"
$ws.Range("E2").Value = "Redux"
$ws.Range("F2").Value = "Router"
$ws.Range("H2").Value = "Knowledge"
$ws.Range("I2").Value = "Notes"

$ws.Range("A3").Value = "Seed"
$ws.Range("B3").Value = "Product, Categories"
$ws.Range("H3").Value = "Sheetjs"

$ws.Range("A4").Value = "Components"

$ws.Range("A5").Value = "  Inventory View All"
$ws.Range("B5").Value = "Product/Categories"
$ws.Range("C5").Value = "GetAll"
$ws.Range("E5").Value = "Inventory"

$ws.Range("A6").Value = "  Inventory View Certain Types"
$ws.Range("B6").Value = "^"
$ws.Range("C6").Value = "GetAllByType"
$ws.Range("E6").Value = "Inventory"

$ws.Range("A7").Value = "Dashboard"
$ws.Range("H7").Value = "D3js/Sheetjs"

$ws.Range("A8").Value = "  Overview"
$ws.Range("B8").Value = "GeneralLedger"

$ws.Range("A9").Value = "  Inventory"
$ws.Range("B9").Value = "Product/Categories"
$ws.Range("C9").Value = "GetAll, GetAllByType"
$ws.Range("I9").Value = "Inventory Turnover (Sales / Inventory), Inventory Write-Off (Discontinued / Inventory"

$ws.Range("A11").Value = "Utilities"

$ws.Range("A12").Value = "  Date from Excel to Javascript"
$ws.Range("I12").Value = "See #1"

$ws.Range("A15").Value = "Tier 2"

$ws.Range("A16").Value = "Components"

$ws.Range("A17").Value = "  Orders View"
$ws.Range("B17").Value = "Customers/Orders/OrderDetails"
$ws.Range("C17").Value = "GetAll"
$ws.Range("E17").Value = "Orders"
$ws.Range("I17").Value = "Ratios - Segment By Demand"

$ws.Range("A18").Value = "  Orders By Customers View"
$ws.Range("B18").Value = "^"
$ws.Range("C18").Value = "GetAllByCustomers"
$ws.Range("E18").Value = "Orders"

$ws.Range("A19").Value = "Dashboard"

$ws.Range("A20").Value = "  Customer"
$ws.Range("B20").Value = "Customers/Orders/OrderDetails"
$ws.Range("C20").Value = "GetAllByCustomers"
$ws.Range("E20").Value = "Orders"
$ws.Range("I20").Value = "Segmentation, Demand, Gross Profit"

$ws.Range("A21").Value = "  Overview"
$ws.Range("I21").Value = "GrossProfit, NetProfit"

$ws.Range("A23").Value = "Tier 3"

$ws.Range("A24").Value = "  Payment System"

$ws.Range("A26").Value = "Ability to make it into mobile app"

$ws.Range("A29").Value = "#1"

$ws.Range("A30").Value = "function getJsDateFromExcel(excelDate) {"


$ws.Range("B32").Value = "// JavaScript dates can be constructed by passing milliseconds"

$ws.Range("B33").Value = "// since the Unix epoch (January 1, 1970) example: new Date(12312512312);"


$ws.Range("B35").Value = "// 1. Subtract number of days between Jan 1, 1900 and Jan 1, 1970, plus 1 (Google `"excel leap year bug`") "

$ws.Range("B36").Value = "// 2. Convert to milliseconds."


$ws.Range("B38").Value = "return new Date((excelDate - (25567 + 1))*86400*1000);"


$ws.Range("A41").Value = "Notes"

$ws.Range("A42").Value = "`"Tech is about pushing yourself`""

# --- Formatting ---

# Highlighted fields completed for the new "Tier 1 Backend" rows (fill = Green,
# Accent 6, Lighter 40% in the source workbook)
$highlightCells = @("B3","B5","C5","B6","C6","B9")
foreach ($addr in $highlightCells) {
    $ws.Range($addr).Interior.ThemeColor = 10
}

# Code-sample block (re-applies the Consolas-based syntax-highlight colors)
$codeKeyword = @("A30","B38")
foreach ($addr in $codeKeyword) {
    $ws.Range($addr).Font.Name = "Consolas"
    $ws.Range($addr).Font.Color = 12665455
}

$codeText = @("A31","B31","A32","A33","A34","B34","A35","A36","A37","B37","A38","A39")
foreach ($addr in $codeText) {
    $ws.Range($addr).Font.Name = "Consolas"
    $ws.Range($addr).Font.Color = 3025188
}

$codeComment = @("B32","B33","B35","B36")
foreach ($addr in $codeComment) {
    $ws.Range($addr).Font.Name = "Consolas"
    $ws.Range($addr).Font.Color = 8221546
}

# --- Sheet view / selection state ---
$ws.Range("C9").Select()
